$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current (pre-edit) row contents for every affected row (columns A:AY)
# before any writes happen, so cycle-based permutations read the original values.
$snapshot = @{}
$snapshot[18] = $ws.Range("A18:AY18").Value2
$snapshot[19] = $ws.Range("A19:AY19").Value2
$snapshot[21] = $ws.Range("A21:AY21").Value2
$snapshot[22] = $ws.Range("A22:AY22").Value2
$snapshot[23] = $ws.Range("A23:AY23").Value2
$snapshot[24] = $ws.Range("A24:AY24").Value2
$snapshot[25] = $ws.Range("A25:AY25").Value2
$snapshot[26] = $ws.Range("A26:AY26").Value2
$snapshot[27] = $ws.Range("A27:AY27").Value2
$snapshot[31] = $ws.Range("A31:AY31").Value2
$snapshot[33] = $ws.Range("A33:AY33").Value2
$snapshot[36] = $ws.Range("A36:AY36").Value2
$snapshot[37] = $ws.Range("A37:AY37").Value2
$snapshot[40] = $ws.Range("A40:AY40").Value2
$snapshot[41] = $ws.Range("A41:AY41").Value2
$snapshot[43] = $ws.Range("A43:AY43").Value2
$snapshot[44] = $ws.Range("A44:AY44").Value2
$snapshot[45] = $ws.Range("A45:AY45").Value2
$snapshot[47] = $ws.Range("A47:AY47").Value2
$snapshot[49] = $ws.Range("A49:AY49").Value2
$snapshot[50] = $ws.Range("A50:AY50").Value2
$snapshot[52] = $ws.Range("A52:AY52").Value2
$snapshot[53] = $ws.Range("A53:AY53").Value2
$snapshot[54] = $ws.Range("A54:AY54").Value2
$snapshot[56] = $ws.Range("A56:AY56").Value2
$snapshot[57] = $ws.Range("A57:AY57").Value2
$snapshot[58] = $ws.Range("A58:AY58").Value2
$snapshot[60] = $ws.Range("A60:AY60").Value2
$snapshot[61] = $ws.Range("A61:AY61").Value2
$snapshot[62] = $ws.Range("A62:AY62").Value2
$snapshot[63] = $ws.Range("A63:AY63").Value2
$snapshot[66] = $ws.Range("A66:AY66").Value2
$snapshot[67] = $ws.Range("A67:AY67").Value2
$snapshot[68] = $ws.Range("A68:AY68").Value2
$snapshot[75] = $ws.Range("A75:AY75").Value2
$snapshot[76] = $ws.Range("A76:AY76").Value2
$snapshot[77] = $ws.Range("A77:AY77").Value2

# Write each destination row with the snapshot of its mapped source row.
$ws.Range("A18:AY18").Value2 = $snapshot[19]
$ws.Range("A19:AY19").Value2 = $snapshot[18]
$ws.Range("A21:AY21").Value2 = $snapshot[26]
$ws.Range("A22:AY22").Value2 = $snapshot[24]
$ws.Range("A23:AY23").Value2 = $snapshot[25]
$ws.Range("A24:AY24").Value2 = $snapshot[23]
$ws.Range("A25:AY25").Value2 = $snapshot[22]
$ws.Range("A26:AY26").Value2 = $snapshot[27]
$ws.Range("A27:AY27").Value2 = $snapshot[21]
$ws.Range("A31:AY31").Value2 = $snapshot[33]
$ws.Range("A33:AY33").Value2 = $snapshot[31]
$ws.Range("A36:AY36").Value2 = $snapshot[37]
$ws.Range("A37:AY37").Value2 = $snapshot[36]
$ws.Range("A40:AY40").Value2 = $snapshot[41]
$ws.Range("A41:AY41").Value2 = $snapshot[40]
$ws.Range("A43:AY43").Value2 = $snapshot[45]
$ws.Range("A44:AY44").Value2 = $snapshot[43]
$ws.Range("A45:AY45").Value2 = $snapshot[44]
$ws.Range("A47:AY47").Value2 = $snapshot[50]
$ws.Range("A49:AY49").Value2 = $snapshot[47]
$ws.Range("A50:AY50").Value2 = $snapshot[49]
$ws.Range("A52:AY52").Value2 = $snapshot[57]
$ws.Range("A53:AY53").Value2 = $snapshot[54]
$ws.Range("A54:AY54").Value2 = $snapshot[52]
$ws.Range("A56:AY56").Value2 = $snapshot[53]
$ws.Range("A57:AY57").Value2 = $snapshot[56]
$ws.Range("A58:AY58").Value2 = $snapshot[61]
$ws.Range("A60:AY60").Value2 = $snapshot[58]
$ws.Range("A61:AY61").Value2 = $snapshot[60]
$ws.Range("A62:AY62").Value2 = $snapshot[63]
$ws.Range("A63:AY63").Value2 = $snapshot[62]
$ws.Range("A66:AY66").Value2 = $snapshot[68]
$ws.Range("A67:AY67").Value2 = $snapshot[66]
$ws.Range("A68:AY68").Value2 = $snapshot[67]
$ws.Range("A75:AY75").Value2 = $snapshot[77]
$ws.Range("A76:AY76").Value2 = $snapshot[75]
$ws.Range("A77:AY77").Value2 = $snapshot[76]

Write-Output "Row content permutation applied."